# "Just say no to suet"
# The ingredient "suet" is removed from the Ingredients sheet. This is a
# plain row deletion: the row holding "suet" (row 161) is removed, every
# row below it shifts up by one, and Excel automatically drops the now
# unused "suet" shared string and renumbers the shared-string table when
# the workbook is saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ingredients")

# Locate the row that holds "suet" in column A (defensive: search instead
# of hard-coding the row number, in case the sheet ever shifts).
$suetRow = 0
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp
for ($r = 1; $r -le $lastRow; $r++) {
    if ($ws.Cells.Item($r, 1).Text -eq "suet") {
        $suetRow = $r
        break
    }
}

if ($suetRow -eq 0) {
    throw "Could not find a row containing 'suet' in column A"
}

# Delete the entire row; remaining rows shift up automatically.
$ws.Rows.Item($suetRow).Delete()

# Keep the AutoFilter's defined name range in sync with the new last row
# (the table lost one data row, e.g. $A$2:$G$178 -> $A$2:$G$177).
$newLastRow = $lastRow - 1
$filterName = $wb.Names.Item("_xlnm._FilterDatabase")
$filterName.RefersTo = "=Ingredients!`$A`$2:`$G`$" + $newLastRow

# Mirror the author's resulting selection: the row that now occupies the
# position the deleted row used to sit in.
$ws.Rows.Item($suetRow).Select()
